# Generate Report for Handback
# The CI report is regenerated: the b1e3693c... file moved from
# "Ready for handoff" to "Handed back: in sync with en-US" (with a new
# handback timestamp), and the two rows (b1e3693c / cbe7af25) are
# re-sorted alphabetically (b1e3693c now first) on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "b1e3693c-9008-434f-9246-cd16926eb1aa.md"
$ov.Range("B2").Value = "Handed back: in sync with en-US"
$ov.Range("C2").Value = "Handed back: in sync with en-US"
$ov.Range("D2").Value = "2016-03-23 08:53:27"

$ov.Range("A3").Value = "cbe7af25-687c-4191-a915-7ed22c306881.md"
$ov.Range("B3").Value = "Handed back: in sync with en-US"
$ov.Range("C3").Value = "Handed back: in sync with en-US"
$ov.Range("D3").Value = "2016-03-23 08:52:00"

$ovLinks = @($ov.Hyperlinks)
foreach ($h in $ovLinks) {
  if ($h.Range.Row -eq 2) { $h.TextToDisplay = "b1e3693c-9008-434f-9246-cd16926eb1aa.md" }
  if ($h.Range.Row -eq 3) { $h.TextToDisplay = "cbe7af25-687c-4191-a915-7ed22c306881.md" }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "b1e3693c-9008-434f-9246-cd16926eb1aa.md"
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("D2").Value = "b1e3693c-9008-434f-9246-cd16926eb1aa.263448d73583d788a29ab2ebfc86ba38fb7ef971.zh-cn.xlf"
$zh.Range("E2").Value = "2016-03-23 08:53:23"
$zh.Range("F2").Value = "b1e3693c-9008-434f-9246-cd16926eb1aa.md"
$zh.Range("G2").Value = "b1e3693c-9008-434f-9246-cd16926eb1aa.263448d73583d788a29ab2ebfc86ba38fb7ef971.zh-cn.xlf"
$zh.Range("H2").Value = "2016-03-23 08:53:47"
$zh.Range("J2").Value = "Include"

$zh.Range("A3").Value = "cbe7af25-687c-4191-a915-7ed22c306881.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Handed back: in sync with en-US"
$zh.Range("D3").Value = "cbe7af25-687c-4191-a915-7ed22c306881.838a0d4ed3ff88b67ff2bd764b96a9e7aca28754.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-23 08:51:56"
$zh.Range("F3").Value = "cbe7af25-687c-4191-a915-7ed22c306881.md"
$zh.Range("G3").Value = "cbe7af25-687c-4191-a915-7ed22c306881.838a0d4ed3ff88b67ff2bd764b96a9e7aca28754.zh-cn.xlf"
$zh.Range("H3").Value = "2016-03-23 08:52:32"
$zh.Range("J3").Value = "Include"

$zhLinks = @($zh.Hyperlinks)
foreach ($h in $zhLinks) {
  $col = $h.Range.Column
  $row = $h.Range.Row
  if ($row -eq 2) {
    if ($col -eq 1) { $h.TextToDisplay = "b1e3693c-9008-434f-9246-cd16926eb1aa.md" }
    if ($col -eq 4) { $h.TextToDisplay = "b1e3693c-9008-434f-9246-cd16926eb1aa.263448d73583d788a29ab2ebfc86ba38fb7ef971.zh-cn.xlf" }
    if ($col -eq 6) { $h.TextToDisplay = "b1e3693c-9008-434f-9246-cd16926eb1aa.md" }
    if ($col -eq 7) { $h.TextToDisplay = "b1e3693c-9008-434f-9246-cd16926eb1aa.263448d73583d788a29ab2ebfc86ba38fb7ef971.zh-cn.xlf" }
  }
  if ($row -eq 3) {
    if ($col -eq 1) { $h.TextToDisplay = "cbe7af25-687c-4191-a915-7ed22c306881.md" }
    if ($col -eq 4) { $h.TextToDisplay = "cbe7af25-687c-4191-a915-7ed22c306881.838a0d4ed3ff88b67ff2bd764b96a9e7aca28754.zh-cn.xlf" }
    if ($col -eq 6) { $h.TextToDisplay = "cbe7af25-687c-4191-a915-7ed22c306881.md" }
    if ($col -eq 7) { $h.TextToDisplay = "cbe7af25-687c-4191-a915-7ed22c306881.838a0d4ed3ff88b67ff2bd764b96a9e7aca28754.zh-cn.xlf" }
  }
}

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "b1e3693c-9008-434f-9246-cd16926eb1aa.md"
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("D2").Value = "b1e3693c-9008-434f-9246-cd16926eb1aa.263448d73583d788a29ab2ebfc86ba38fb7ef971.de-de.xlf"
$de.Range("E2").Value = "2016-03-23 08:53:27"
$de.Range("F2").Value = "b1e3693c-9008-434f-9246-cd16926eb1aa.md"
$de.Range("G2").Value = "b1e3693c-9008-434f-9246-cd16926eb1aa.263448d73583d788a29ab2ebfc86ba38fb7ef971.de-de.xlf"
$de.Range("H2").Value = "2016-03-23 08:53:54"
$de.Range("J2").Value = "Include"

$de.Range("A3").Value = "cbe7af25-687c-4191-a915-7ed22c306881.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Handed back: in sync with en-US"
$de.Range("D3").Value = "cbe7af25-687c-4191-a915-7ed22c306881.838a0d4ed3ff88b67ff2bd764b96a9e7aca28754.de-de.xlf"
$de.Range("E3").Value = "2016-03-23 08:52:00"
$de.Range("F3").Value = "cbe7af25-687c-4191-a915-7ed22c306881.md"
$de.Range("G3").Value = "cbe7af25-687c-4191-a915-7ed22c306881.838a0d4ed3ff88b67ff2bd764b96a9e7aca28754.de-de.xlf"
$de.Range("H3").Value = "2016-03-23 08:52:41"
$de.Range("J3").Value = "Include"

$deLinks = @($de.Hyperlinks)
foreach ($h in $deLinks) {
  $col = $h.Range.Column
  $row = $h.Range.Row
  if ($row -eq 2) {
    if ($col -eq 1) { $h.TextToDisplay = "b1e3693c-9008-434f-9246-cd16926eb1aa.md" }
    if ($col -eq 4) { $h.TextToDisplay = "b1e3693c-9008-434f-9246-cd16926eb1aa.263448d73583d788a29ab2ebfc86ba38fb7ef971.de-de.xlf" }
    if ($col -eq 6) { $h.TextToDisplay = "b1e3693c-9008-434f-9246-cd16926eb1aa.md" }
    if ($col -eq 7) { $h.TextToDisplay = "b1e3693c-9008-434f-9246-cd16926eb1aa.263448d73583d788a29ab2ebfc86ba38fb7ef971.de-de.xlf" }
  }
  if ($row -eq 3) {
    if ($col -eq 1) { $h.TextToDisplay = "cbe7af25-687c-4191-a915-7ed22c306881.md" }
    if ($col -eq 4) { $h.TextToDisplay = "cbe7af25-687c-4191-a915-7ed22c306881.838a0d4ed3ff88b67ff2bd764b96a9e7aca28754.de-de.xlf" }
    if ($col -eq 6) { $h.TextToDisplay = "cbe7af25-687c-4191-a915-7ed22c306881.md" }
    if ($col -eq 7) { $h.TextToDisplay = "cbe7af25-687c-4191-a915-7ed22c306881.838a0d4ed3ff88b67ff2bd764b96a9e7aca28754.de-de.xlf" }
  }
}
